$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Handback report generation:
#   - Overview / per-locale "Status" columns flip from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - zh-cn / de-de rows gain their "Latest Target File" (hyperlinked source
#     doc), "Latest Handback File" (xlf) and "Latest Handback DateTime" values
#   - columns that now hold longer strings are widened to fit
# ---------------------------------------------------------------------------

$statusText = "Handed back: in sync with en-US"
$mdFileName = "12ff1f41-8210-45e4-bddd-5ef673472969.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af61115601566d82751cf29e7b82e1dacac981a2/e2e/12ff1f41-8210-45e4-bddd-5ef673472969.md"

# -- Overview sheet ----------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.16
$wsOverview.Columns.Item(6).ColumnWidth = 29.16

# -- zh-cn sheet ---------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("I2").Value = $mdFileName
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl, "", "", $mdFileName)
$wsZh.Range("J2").Value = "12ff1f41-8210-45e4-bddd-5ef673472969.03874216514a5b8b216f9db17b7bd6044f37df7f.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-26 09:04:31"

$wsZh.Columns.Item(3).ColumnWidth = 29.16
$wsZh.Columns.Item(9).ColumnWidth = 39.16
$wsZh.Columns.Item(10).ColumnWidth = 39.16

# -- de-de sheet ---------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("I2").Value = $mdFileName
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl, "", "", $mdFileName)
$wsDe.Range("J2").Value = "12ff1f41-8210-45e4-bddd-5ef673472969.03874216514a5b8b216f9db17b7bd6044f37df7f.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-26 09:04:39"

$wsDe.Columns.Item(3).ColumnWidth = 29.16
$wsDe.Columns.Item(9).ColumnWidth = 39.16
$wsDe.Columns.Item(10).ColumnWidth = 39.16
